$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "color" column with the per-career color codes (C1:C7).
# The previous (empty, time-formatted) C2:C7 cells need their old
# number format cleared so the new text values land with default style.
$ws.Range("C2:C7").ClearFormats()

$ws.Range("C1").Value = "color"
$ws.Range("C2").Value = "2257B9"
$ws.Range("C3").Value = "22B949"
$ws.Range("C4").Value = "22B9B7"
$ws.Range("C5").Value = "22B964"
$ws.Range("C6").Value = "E9A660"
$ws.Range("C7").Value = "5D63EA"

# Leave the selection where the author left it after editing.
$ws.Range("I8").Select()
